$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is numeric-looking need to be forced to Text
# format so Excel stores them as strings (matching the source inline-string
# cells) instead of auto-converting to a number. The style is reset back to
# "Normal" immediately after the write so no stray number-format style sticks
# around on the cell (matches the target: no style/formatting changes).

$ws.Range("D2").Value = '27.131.64'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '1.898.96'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5229'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3805'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07285'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.36'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9025'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08172'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = '1.851.74'
$ws.Range("E14").Value = '  -2.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.353'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008645'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = '27.169.36'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.116'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.468'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.327'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.744'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.829'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.897'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09208'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05040'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7949'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.220'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.956'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.365'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.643'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5703'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.082'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.065'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.588'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("E51").Value = '  +0.47%  '
